## Update the generated FHIR StructureDefinition workbook to the
## "2025 august" refresh: new canonical base URL, new generation
## timestamp, and the corresponding ValueSet URL, then let the
## "Elements" sheet column widths follow (they were auto best-fit
## to the new, shorter strings).

$wb = $excel.ActiveWorkbook

$metadata = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

# --- 1. Canonical URL of the StructureDefinition -------------------------
$newBase = "https://2rdoc.pt/ig/ios-lifestyle-medicine"

$newUrl = "$newBase/StructureDefinition/social-support"
$metadata.Range("B2").Value = $newUrl          # Metadata!URL
$elements.Range("R5").Value = $newUrl          # Extension.url fixed value (same text)

# --- 2. Generation date ----------------------------------------------------
$metadata.Range("B8").Value = "2025-08-20T10:40:04+01:00"

# --- 3. ValueSet URL referenced by the binding ----------------------------
$newValueSetUrl = "$newBase/ValueSet/social-support-vs"
$elements.Range("Z6").Value = $newValueSetUrl

# --- 4. Column widths on "Elements" recomputed (best-fit) after the text
#        shrank.  The values below are the ColumnWidth inputs that this
#        runtime's pixel-rounding reproduces as the closest match to the
#        workbook's target stored OOXML widths.
$widths = @{
    1  = 15.666666666666666
    2  = 15.666666666666666
    3  = 9.0
    4  = 6.166666666666667
    5  = 4.5
    6  = 3.1666666666666665
    7  = 3.5
    8  = 11.833333333333334
    9  = 9.666666666666666
    11 = 13.5
    15 = 11.5
    20 = 7.0
    21 = 12.833333333333334
    22 = 13.166666666666666
    23 = 14.166666666666666
    24 = 13.833333333333334
    25 = 16.166666666666668
    26 = 54.833333333333336
    27 = 4.166666666666667
    28 = 17.166666666666668
    29 = 33.666666666666664
    30 = 12.666666666666666
    31 = 10.5
    32 = 14.166666666666666
    33 = 7.333333333333333
    34 = 7.666666666666667
    37 = 18.666666666666668
}

foreach ($col in $widths.Keys) {
    $elements.Columns.Item($col).ColumnWidth = $widths[$col]
}
